$wb = $excel.ActiveWorkbook

# --- Sheet "c meal mixer" (sheet1): selection moves from G18 to A6 ---
$wsMeal = $wb.Worksheets.Item(1)
$wsMeal.Range("A6").Select()

# --- Sheet "c clinker kiln" (sheet2): add "slag in kiln" calculations ---
$wsKiln = $wb.Worksheets.Item(2)

# New header columns
$wsKiln.Range("G1").Value = "2nd Known Substance"
$wsKiln.Range("H1").Value = "2Qty Origin"

# Row 9: clinker output -> clinker_co2 (Addition with CO2 output)
$wsKiln.Range("A9").Value = "clinker"
$wsKiln.Range("C9").Value = "clinker_co2"
$wsKiln.Range("E9").Value = "Addition"
$wsKiln.Range("F9").Value = "NONE"
$wsKiln.Range("B9").Value = "output"
$wsKiln.Range("D9").Value = "tmp"
$wsKiln.Range("G9").Value = "CO2"
$wsKiln.Range("H9").Value = "output"

# Row 10: meal inflow -> slag (Subtraction with clinker_co2 tmp)
$wsKiln.Range("A10").Value = "meal"
$wsKiln.Range("C10").Value = "slag"
$wsKiln.Range("E10").Value = "Subtraction"
$wsKiln.Range("F10").Value = "NONE"
$wsKiln.Range("B10").Value = "inflow"
$wsKiln.Range("D10").Value = "output"
$wsKiln.Range("G10").Value = "clinker_co2"
$wsKiln.Range("H10").Value = "tmp"

# --- Sheet "var clinker kiln" (sheet5): selection moves from F7 to E6 ---
$wsKilnVar = $wb.Worksheets.Item(5)
$wsKilnVar.Range("E6").Select()

# --- Finally activate "c clinker kiln" and select I14 (becomes the active tab) ---
$wsKiln.Activate()
$wsKiln.Range("I14").Select()
